$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 66, pushing the existing rows 66-116 down to 67-117
$ws.Rows.Item(66).EntireRow.Insert()

# Populate the newly inserted row 66 with the new weekly record
$ws.Cells.Item(66, 1).Value = 1
$ws.Cells.Item(66, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(66, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(66, 4).Value = 44904
$ws.Cells.Item(66, 5).Value = 15
$ws.Cells.Item(66, 6).Value = 100112021
$ws.Cells.Item(66, 7).Value = "Ají"
$ws.Cells.Item(66, 8).Value = "Inferno"
$ws.Cells.Item(66, 9).Value = "Primera"
$ws.Cells.Item(66, 10).Value = 140
$ws.Cells.Item(66, 11).Value = 9000
$ws.Cells.Item(66, 12).Value = 10000
$ws.Cells.Item(66, 13).Value = 9500
$ws.Cells.Item(66, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(66, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(66, 16).Value = 633
$ws.Cells.Item(66, 17).Value = 15
$ws.Cells.Item(66, 18).Value = "Hortaliza"
